$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53:73 down to 54:74.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with this week's data point.
$ws.Range("A53").Value = 10
$ws.Range("B53").Value = "Vega Modelo de Temuco"
$ws.Range("C53").Value = "La Araucanía"
$ws.Range("D53").Value = 44609
$ws.Range("E53").Value = 9
$ws.Range("F53").Value = 100112030
$ws.Range("G53").Value = "Poroto granado"
$ws.Range("H53").Value = "Sin especificar"
$ws.Range("I53").Value = "Primera"
$ws.Range("J53").Value = 10
$ws.Range("K53").Value = 25000
$ws.Range("L53").Value = 25000
$ws.Range("M53").Value = 25000
$ws.Range("N53").Value = "$/saco 25 kilos"
$ws.Range("O53").Value = "Región de La Araucanía"
$ws.Range("P53").Value = 1000
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = "Hortaliza"
